# Update footer version number and last-update date on the Quick Reference
# footer (footer2.xml): "Version 3.0.x" -> "Version 3.1.x"
# and the cached DATE field result "2024-07-02" -> "2024-09-18".

$d = $word.ActiveDocument

# 1) Bump the minor version digit: "3.0.x" -> "3.1.x"
$d.Content.Find.Execute("3.0.x", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.1.x", 2)

# 2) Update the "Last update" date text (cached field result)
$d.Content.Find.Execute("2024-07-02", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-18", 2)
